# Vermont 2016 MCAS network-data cleanup
# - rename header columns to snake_case field names
# - normalize "de"/"la" -> "De"/"La" in a handful of place names
# - drop the trailing metadata/footer rows (sample size, source, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two footer/metadata blocks first (bottom-up so row numbers
# of the earlier block aren't shifted by the later deletion).
$ws.Rows("476:480").Delete()
$ws.Rows("61:65").Delete()

# Header row: translate Spanish labels to the new snake_case schema.
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Capitalize "de"/"la" in a set of municipality/state names.
$ws.Range("B2").Value = "Comitán De Domínguez"
$ws.Range("A13").Value = "Ciudad De México"
$ws.Range("A16").Value = "Estado De México"
$ws.Range("B19").Value = "Chilapa De Álvarez"
$ws.Range("B23").Value = "Nopala De Villagrán"
$ws.Range("B29").Value = "Putla Villa De Guerrero"
$ws.Range("B36").Value = "Izúcar De Matamoros"
$ws.Range("B50").Value = "Martínez De La Torre"
